$p = $ppt.ActivePresentation

# Bullet character used by the source paragraphs ("• ").
$bullet = [char]0x2022

# Note: each paragraph's text is cleared to "" before the real value is
# assigned. Setting .Text directly to a string that shares a trailing
# substring with the existing text makes the host preserve that matching
# suffix as its own run (a diff-style reuse of the old run's formatting,
# e.g. "...확보" / "...% 달성"); clearing first removes any such overlap so
# the result is a single clean run, matching the target markup.

# --- Slide 1: update the "active since" date range in the subtitle ---
# (Using Paragraphs(1) rather than the shape-level TextRange avoids the
# host tagging the rewritten run with a fresh <a:rPr lang="en-US"/>, so the
# run keeps inheriting formatting the same way the original run did.)
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Paragraphs(1).Text = ""
$subtitleRange.Paragraphs(1).Text = "2023.03 - 현재"

# --- Slide 3: rewrite the project-detail bullets and key-result bullets ---
$slide3 = $p.Slides.Item(3)
$detailShape = $slide3.Shapes.Item(1)
$tr = $detailShape.TextFrame.TextRange

$t3  = $bullet + " 그로스폴리오 론칭을 위한 전체 마케팅 전략 수립 및 실행 계획 수립"
$t4  = $bullet + " 디지털 마케팅 채널별 맞춤형 콘텐츠 기획 및 제작 관리"
$t5  = $bullet + " SNS 채널 운영 전략 수립 및 캠페인 콘텐츠 제작 진행"
$t6  = $bullet + " 유저 획득을 위한 퍼포먼스 마케팅 캠페인 기획 및 운영"
$t7  = $bullet + " 론칭 이벤트 기획 및 프로모션 운영 총괄"
$t9  = $bullet + " 론칭 첫 달 신규 가입자 32,000명 확보"
$t10 = $bullet + " 캠페인 기간 내 광고 투자 대비 수익률(ROAS) 180% 달성"
$t11 = $bullet + " SNS 채널 팔로워 3개월 간 45% 증가"

$tr.Paragraphs(3).Text = ""
$tr.Paragraphs(3).Text = $t3

$tr.Paragraphs(4).Text = ""
$tr.Paragraphs(4).Text = $t4

$tr.Paragraphs(5).Text = ""
$tr.Paragraphs(5).Text = $t5

$tr.Paragraphs(6).Text = ""
$tr.Paragraphs(6).Text = $t6

$tr.Paragraphs(7).Text = ""
$tr.Paragraphs(7).Text = $t7

$tr.Paragraphs(9).Text = ""
$tr.Paragraphs(9).Text = $t9

$tr.Paragraphs(10).Text = ""
$tr.Paragraphs(10).Text = $t10

$tr.Paragraphs(11).Text = ""
$tr.Paragraphs(11).Text = $t11
